# Fruta / hortaliza, semanal
# Weekly data refresh: insert a new observation row at row 25 (pushing the
# existing rows 25-59 down to 26-60) and populate it with this week's values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row at position 25; rows 25..59 shift down to 26..60.
$ws.Rows.Item(25).Insert()

# Populate the newly inserted row 25 with the new weekly record.
$ws.Range("A25").Value = 8
$ws.Range("B25").Value = "Terminal La Palmera de La Serena"
$ws.Range("C25").Value = "Coquimbo"
$ws.Range("D25").Value = 44799
$ws.Range("E25").Value = 4
$ws.Range("F25").Value = 100114007
$ws.Range("G25").Value = "Jengibre"
$ws.Range("H25").Value = "Sin especificar"
$ws.Range("I25").Value = "Primera"
$ws.Range("J25").Value = 460
$ws.Range("K25").Value = 14000
$ws.Range("L25").Value = 15000
$ws.Range("M25").Value = 14500
$ws.Range("N25").Value = "$/caja 13 kilos"
$ws.Range("O25").Value = "Perú"
$ws.Range("P25").Value = 1115
$ws.Range("Q25").Value = 13
$ws.Range("R25").Value = "Hortaliza"
